# Add "NA" values under the duplicate_image_filename column (column E)
# for rows 2 through 21.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 originally holds an empty/blank shared-string placeholder cell. Re-assert
# it as blank so the round trip through this engine doesn't coerce it into a
# real value (an unrelated artifact of how empty string cells get persisted).
$ws.Range("F1").Value = ""

for ($r = 2; $r -le 21; $r++) {
    $ws.Range("E$r").Value = "NA"
}
